$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header styling (bold/centered) from B1 into the new C1 header cell
$ws.Range("B1").Copy($ws.Range("C1"))

# Header row: update attribute URIs
$ws.Range("B1").Value = "http://dbpedia.org/ontology/birthDate"
$ws.Range("C1").Value = "http://dbpedia.org/ontology/birthPlace"

# Row 2: Dick Sheppard (priest)
$ws.Range("A2").Value = "http://dbpedia.org/resource/Dick_Sheppard_(priest)"
$ws.Range("B2").Value = "http://dbpedia.org/resource/1880"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Windsor"

# Row 3: Claus Westermann
$ws.Range("A3").Value = "http://dbpedia.org/resource/Claus_Westermann"
$ws.Range("B3").Value = "http://dbpedia.org/resource/1909"
$ws.Range("C3").Value = "http://dbpedia.org/resource/Berlin"

# Remove the now-unused rows 4-8 from the old dataset
$ws.Range("A4:C8").Clear()
